# Power_FACTS.xlsx - "Fix pyomo vLineP, vLineQ, vGenQ handeling. Add bounds
# for quadratic variables"
#
# The example input sheet lists one row per STATCOM (FACTS) device with an
# "EnableInvest" flag ([0,1]) in column H. All nine STATCOM units
# (STATCOM_1 .. STATCOM_9, sheet rows 7-15) had that flag hard-coded to 1
# (investment enabled) - flip them to 0 so the example dataset matches the
# fixed Pyomo model's expected bounded/quadratic-variable handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power FACTS")

# EnableInvest column (H) for STATCOM_1..STATCOM_9 -> 0 (rows 7 through 15)
$ws.Range("H7:H15").Value = 0

# Leave the sheet selection where the author last left it before saving.
$ws.Range("M24").Select()
